$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recurring Shifts")

function Frac($h, $m) { return ($h * 60 + $m) / 1440 }

# Row 2: Monday 12:00 - 13:00, Weekly Meeting
$ws.Range("A2").Value = "Monday"
$ws.Range("B2").Value = Frac 12 0
$ws.Range("C2").Value = Frac 13 0
$ws.Range("E2").Value = "Weekly Meeting"

# Row 3: Monday 13:30 - 16:00, OH
$ws.Range("A3").Value = "Monday"
$ws.Range("B3").Value = Frac 13 30
$ws.Range("C3").Value = Frac 16 0
$ws.Range("E3").Value = "OH"

# Row 4: Tuesday 10:30 - 13:00, OH
$ws.Range("A4").Value = "Tuesday"
$ws.Range("B4").Value = Frac 10 30
$ws.Range("C4").Value = Frac 13 0
$ws.Range("E4").Value = "OH"

# Row 5: Wednesday 12:00 - 14:30, OH
$ws.Range("A5").Value = "Wednesday"
$ws.Range("B5").Value = Frac 12 0
$ws.Range("C5").Value = Frac 14 30
$ws.Range("E5").Value = "OH"

# Row 6: Thursday 11:30 - 13:00, OH
$ws.Range("A6").Value = "Thursday"
$ws.Range("B6").Value = Frac 11 30
$ws.Range("C6").Value = Frac 13 0
$ws.Range("E6").Value = "OH"

# Row 7 (new): Thursday 15:00 - 16:00, OH
$ws.Range("A7").Value = "Thursday"
$ws.Range("B7").Value = Frac 15 0
$ws.Range("C7").Value = Frac 16 0
$ws.Range("D7").Value = "DesignHub"
$ws.Range("E7").Value = "OH"

# Row 8 (new): Friday 11:00 - 15:00, OH
$ws.Range("A8").Value = "Friday"
$ws.Range("B8").Value = Frac 11 0
$ws.Range("C8").Value = Frac 15 0
$ws.Range("D8").Value = "DesignHub"
$ws.Range("E8").Value = "OH"

# Rows 3 and 4 previously had stale/copy-pasted (C2-B2) and (C3-B3) formulas;
# fix them up so each row computes off its own Clock-in/Clock-out pair.
$ws.Range("G3").Formula = "=(C3-B3)*24"
$ws.Range("G4").Formula = "=(C4-B4)*24"

# Update the Total Hours sum range to include the two new rows
$ws.Range("I3").Formula = "=SUM(G2:G10)"

# Update selection to D10
$ws.Range("D10").Select()

$wb.Save()
